$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove Wu, Qirui's row entirely (row 10). This shifts Zlitni, Hanane
# (old row 11) up to row 10, and the two trailing blank rows (old 12-13)
# up to rows 11-12.
$ws.Rows("10:10").Delete()

# --- Re-point the diagonal "reviewer rotation" formulas in F:J for every
# data row (3-9) so each one continues to refer to the next row's A/F/G/H/I
# cell, now that the row count dropped by one. Row 10 wraps back to row 3.
for ($r = 3; $r -le 9; $r++) {
  $next = $r + 1
  $ws.Range("F$r").Formula = "=A$next"
  $ws.Range("G$r").Formula = "=F$next"
  $ws.Range("H$r").Formula = "=G$next"
  $ws.Range("I$r").Formula = "=H$next"
  $ws.Range("J$r").Formula = "=I$next"
}
$ws.Range("F10").Formula = "=A3"
$ws.Range("G10").Formula = "=F3"
$ws.Range("H10").Formula = "=G3"
$ws.Range("I10").Formula = "=H3"
$ws.Range("J10").Formula = "=I3"

# --- Assign Sekis, Karol (row 7) a repo link now that they're a reviewer
# for the written docs.
$ws.Range("B7").Value = "https://github.com/karolserkis/CAS-741-Pendula"

# --- Print area no longer needs to include the now-removed row.
$ws.PageSetup.PrintArea = '$A$1:$B$10'

# --- Match the author's final selection / window position.
$ws.Range("A10").Select()
$excel.ActiveWindow.Left = 1340
